$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Old Fleet Number"
$ws.Range("N1").Value = "Old Rego"
$ws.Range("O1").Value = "Registered Owner"
$ws.Range("P1").Value = "Roadside Assistance"
$ws.Range("Q1").Value = "PDD number"
